# Updates the "Pais" COVID-19 summary sheet with a newer data pull:
#   - refreshes the "last updated" timestamp in A1
#   - re-sorts / refreshes the per-country rows (A4:H219) by Casos totales desc,
#     which reorders Venezuela/Irlanda, Paraguay/Luxemburgo, Polinesia Francesa/
#     Butan/Antigua y Barbuda/Liechtenstein/San Martin (Parte Francesa), Timor
#     Oriental/Santa Lucia and Montserrat/Islas Malvinas, and updates the
#     underlying case counts for Alemania, Kazajistan, Belgica, Honduras,
#     Venezuela/Irlanda, Paraguay/Luxemburgo, Eritrea, Polinesia
#     Francesa/Butan/Antigua y Barbuda/Liechtenstein/San Martin (Parte
#     Francesa) and Montserrat/Islas Malvinas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 11 de Agosto de 2020 a las 05:15"

# Columns: Pais, Casos totales, Nuevos casos, Casos activos, Recuperados,
#          Casos criticos, Muertes hoy, Muertes
$data = @(
    @('Estados Unidos',5251446,0,2715934,2369320,0,0,166192),
    @('Brasil',3057470,0,2163812,791801,0,0,101857),
    @('India',2267153,0,1581640,640160,0,0,45353),
    @('Rusia',892654,0,696681,180972,0,0,15001),
    @('Sudafrica',563598,0,417200,135777,0,0,10621),
    @('Mexico',485836,5558,327993,104840,0,705,53003),
    @('Peru',478024,0,324020,132932,0,0,21072),
    @('Colombia',397623,0,221485,162984,0,0,13154),
    @('Chile',375044,0,347342,17563,0,0,10139),
    @('España',370060,0,0,0,0,0,28576),
    @('Iran',328844,0,286642,23586,0,0,18616),
    @('Reino Unido',311641,0,0,0,0,0,46526),
    @('Arabia Saudita',289947,0,253478,33270,0,0,3199),
    @('Pakistan',284660,0,260764,17799,0,0,6097),
    @('Banglades',260507,0,150437,106632,0,0,3438),
    @('Argentina',253868,0,108242,140862,0,0,4764),
    @('Italia',250825,0,202248,13368,0,0,35209),
    @('Turquia',241997,0,224970,11169,0,0,5858),
    @('Alemania',218500,0,198900,10335,0,0,9265),
    @('Francia',202775,0,82836,89599,0,0,30340),
    @('Irak',153599,0,109790,38345,0,0,5464),
    @('Filipinas',136638,0,68159,66185,0,0,2294),
    @('Indonesia',127083,0,82236,39082,0,0,5765),
    @('Canada',120132,0,106355,4790,0,0,8987),
    @('Catar',113262,0,109993,3081,0,0,188),
    @('Kazajistan',100164,722,73702,25404,0,0,1058),
    @('Egipto',95666,0,53779,36852,0,0,5035),
    @('Ecuador',94701,0,78608,10161,0,0,5932),
    @('Bolivia',91635,1636,30823,57100,0,72,3712),
    @('Israel',84722,0,58998,25111,0,0,613),
    @('China',84712,44,79284,794,0,0,4634),
    @('Suecia',82972,0,0,0,0,0,5766),
    @('Ucrania',81957,0,44359,35676,0,0,1922),
    @('Oman',81787,0,76124,5142,0,0,521),
    @('Republica Dominicana',80499,0,44910,34261,0,0,1328),
    @('Panama',75394,0,49510,24220,0,0,1664),
    @('Belgica',74620,468,17792,46949,0,7,9879),
    @('Kuwait',72400,0,64028,7890,0,0,482),
    @('Bielorrusia',68947,0,64991,3367,0,0,589),
    @('Emiratos Arabes Unidos',62704,0,56766,5581,0,0,357),
    @('Rumania',62547,0,30311,29507,0,0,2729),
    @('Paises Bajos',59194,0,0,0,0,0,6157),
    @('Guatemala',56987,0,45589,9176,0,0,2222),
    @('Singapur',55292,0,49609,5656,0,0,27),
    @('Portugal',52825,0,38600,12466,0,0,1759),
    @('Polonia',52410,0,36877,13724,0,0,1809),
    @('Japon',47990,0,33058,13885,0,0,1047),
    @('Honduras',47872,418,6649,39717,0,11,1506),
    @('Nigeria',46867,0,33346,12571,0,0,950),
    @('Barein',44397,0,41209,3025,0,0,163),
    @('Ghana',41212,0,38727,2270,0,0,215),
    @('Armenia',40433,0,32616,7021,0,0,796),
    @('Kirguistan',40085,0,32126,6485,0,0,1474),
    @('Afganistan',37162,0,26228,9606,0,0,1328),
    @('Suiza',36708,0,32400,2321,0,0,1987),
    @('Argelia',35712,0,24920,9480,0,0,1312),
    @('Marruecos',34063,0,24524,9023,0,0,516),
    @('Azerbaiyan',33647,0,30642,2513,0,0,492),
    @('Uzbekistan',31304,0,22992,8112,0,0,200),
    @('Serbia',28262,0,18965,8651,0,0,646),
    @('Moldavia',27841,0,19300,7691,0,0,850),
    @('Kenia',26928,0,13495,13010,0,0,423),
    @('Venezuela',26800,0,16930,9641,0,0,229),
    @('Irlanda',26768,0,23364,1632,0,0,1772),
    @('Costa Rica',23872,0,7823,15805,0,0,244),
    @('Etiopia',23591,0,10411,12760,0,0,420),
    @('Nepal',23310,0,16493,6738,0,0,79),
    @('Austria',22106,0,20010,1373,0,0,723),
    @('Australia',21750,353,12144,9274,0,19,332),
    @('El Salvador',20872,0,9720,10589,0,0,563),
    @('Chequia',18494,0,12983,5121,0,0,390),
    @('Camerun',18042,0,15320,2327,0,0,395),
    @('Costa de Marfil',16798,0,13052,3641,0,0,105),
    @('Dinamarca',14815,0,12925,1270,0,0,620),
    @('Corea del Sur',14660,34,13729,626,0,0,305),
    @('Estado de Palestina',14510,0,8045,6365,0,0,100),
    @('Bosnia y Herzegovina',14498,0,8159,5914,0,0,425),
    @('Bulgaria',13512,0,7980,5073,0,0,459),
    @('Madagascar',13202,0,11011,2040,0,0,151),
    @('Sudan',11956,0,6266,4909,0,0,781),
    @('Republica de Macedonia',11942,0,8087,3327,0,0,528),
    @('Senegal',11312,0,7390,3686,0,0,236),
    @('Noruega',9684,0,8857,571,0,0,256),
    @('Consejo Danes para los Refugiados',9489,0,8363,902,0,0,224),
    @('Malasia',9094,0,8803,166,0,0,125),
    @('Guayana Francesa',8324,0,7549,726,0,0,49),
    @('Zambia',8210,0,6802,1167,0,0,241),
    @('Gabon',8006,0,5823,2132,0,0,51),
    @('Guinea',7930,0,6898,982,0,0,50),
    @('Tayikistan',7827,0,6614,1151,0,0,62),
    @('Haiti',7634,0,4982,2469,0,0,183),
    @('Finlandia',7601,0,6980,288,0,0,333),
    @('Paraguay',7234,0,5276,1876,0,0,82),
    @('Luxemburgo',7216,0,6170,925,0,0,121),
    @('Libano',6812,0,2290,4442,0,0,80),
    @('Mauritania',6555,0,5570,828,0,0,157),
    @('Albania',6536,0,3379,2957,0,0,200),
    @('Libia',5929,0,724,5080,0,0,125),
    @('Grecia',5749,0,3804,1732,0,0,213),
    @('Croacia',5649,0,4906,585,0,0,158),
    @('Republica de Yibuti',5347,0,5120,168,0,0,59),
    @('Maldivas',5157,0,2835,2303,0,0,19),
    @('Guinea Ecuatorial',4821,0,2182,2556,0,0,83),
    @('Zimbabue',4748,0,1524,3120,0,0,104),
    @('Hungria',4731,0,3525,601,0,0,605),
    @('Malaui',4674,0,2430,2098,0,0,146),
    @('Republica de Africa Central',4641,0,1721,2860,0,0,60),
    @('Hong Kong',4149,0,2916,1178,0,0,55),
    @('Nicaragua',3902,0,2913,866,0,0,123),
    @('Montenegro',3696,0,2521,1107,0,0,68),
    @('Congo',3664,0,1589,2017,0,0,58),
    @('Tailandia',3351,0,3160,133,0,0,58),
    @('Suazilandia',3309,0,1634,1614,0,0,61),
    @('Somalia',3227,0,1728,1406,0,0,93),
    @('Namibia',3101,0,715,2367,0,0,19),
    @('Mayotte',3068,0,2835,194,0,0,39),
    @('Cuba',3046,0,2460,498,0,0,88),
    @('Cabo Verde',2883,0,2128,723,0,0,32),
    @('Sri Lanka',2871,0,2593,267,0,0,11),
    @('Eslovaquia',2599,0,1866,702,0,0,31),
    @('Mali',2573,0,1969,479,0,0,125),
    @('Surinam',2489,0,1674,785,0,0,30),
    @('Sudan del Sur',2470,0,1175,1248,0,0,47),
    @('Mozambique',2411,0,860,1535,0,0,16),
    @('Lituania',2265,0,1670,514,0,0,81),
    @('Eslovenia',2255,0,1960,167,0,0,128),
    @('Estonia',2158,0,1962,133,0,0,63),
    @('Ruanda',2152,0,1392,753,0,0,7),
    @('Guinea-Bisau',2052,0,944,1079,0,0,29),
    @('Islandia',1962,0,1838,114,0,0,10),
    @('Benin',1936,0,1600,298,0,0,38),
    @('Sierra Leona',1917,0,1447,401,0,0,69),
    @('Yemen',1832,0,915,399,0,0,518),
    @('Tunez',1717,0,1265,401,0,0,51),
    @('Angola',1679,0,569,1032,0,0,78),
    @('Nueva Zelanda',1570,1,1526,22,0,0,22),
    @('Uruguay',1364,0,1146,181,0,0,37),
    @('Uganda',1297,0,1137,151,0,0,9),
    @('Letonia',1290,0,1070,188,0,0,32),
    @('Jordania',1268,0,1187,70,0,0,11),
    @('Siria',1255,0,364,839,0,0,52),
    @('Republica de Chipre',1252,0,870,363,0,0,19),
    @('Georgia',1250,0,1010,223,0,0,17),
    @('Liberia',1240,0,725,436,0,0,79),
    @('Gambia',1235,0,221,991,0,0,23),
    @('Burkina Faso',1204,0,984,166,0,0,54),
    @('Niger',1158,0,1062,27,0,0,69),
    @('Malta',1112,0,688,415,0,0,9),
    @('Togo',1067,0,729,313,0,0,25),
    @('Botsuana',1066,0,80,984,0,0,2),
    @('Jamaica',1031,8,745,272,0,0,14),
    @('Principado de Andorra',963,0,839,72,0,0,52),
    @('Bahamas',945,0,113,817,0,0,15),
    @('Republica del Chad',945,0,843,26,0,0,76),
    @('Santo Tome y Principe',878,0,800,63,0,0,15),
    @('Vietnam',847,0,399,433,0,0,15),
    @('Lesoto',781,0,175,582,0,0,24),
    @('Crucero',712,0,651,48,0,0,13),
    @('San Marino',699,0,657,0,0,0,42),
    @('Reunion',690,0,631,54,0,0,5),
    @('Aruba',630,0,114,513,0,0,3),
    @('Guyana',568,0,189,357,0,0,22),
    @('Tanzania',509,0,183,305,0,0,21),
    @('Taiwan',480,0,443,30,0,0,7),
    @('Burundi',408,0,315,92,0,0,1),
    @('Comoras',399,0,369,23,0,0,7),
    @('Birmania',360,0,312,42,0,0,6),
    @('Mauricio',344,0,334,0,0,0,10),
    @('Martinica',336,0,98,222,0,0,16),
    @('Isla de Man',336,0,312,0,0,0,24),
    @('Guadalupe',317,0,186,117,0,0,14),
    @('Islas Feroe',306,0,215,91,0,0,0),
    @('Mongolia',293,0,263,30,0,0,0),
    @('Eritrea',285,0,248,37,0,0,0),
    @('Trinidad yTobago',281,0,138,135,0,0,8),
    @('Camboya',251,0,219,32,0,0,0),
    @('Islas Turcas y Caicos',216,0,39,175,0,0,2),
    @('Papua Nueva Guinea',214,0,53,158,0,0,3),
    @('Islas Caimanes',203,0,202,0,0,0,1),
    @('Gibraltar',201,0,186,15,0,0,0),
    @('San Martin (Parte Holandesa)',189,0,86,86,0,0,17),
    @('Bermudas',158,0,144,5,0,0,9),
    @('Belice',154,0,32,120,0,0,2),
    @('Barbados',142,0,112,23,0,0,7),
    @('Brunei',142,0,138,1,0,0,3),
    @('Monaco',133,0,113,16,0,0,4),
    @('Seychelles',126,0,125,1,0,0,0),
    @('Polinesia Francesa',112,0,62,50,0,0,0),
    @('Butan',110,0,96,14,0,0,0),
    @('Antigua y Barbuda',92,0,76,13,0,0,3),
    @('Liechtenstein',89,0,85,3,0,0,1),
    @('San Martin (Parte Francesa)',78,0,44,31,0,0,3),
    @('San Vicente y las Granadinas',57,0,52,5,0,0,0),
    @('Macao',46,0,46,0,0,0,0),
    @('Puerto Rico',39,0,1,36,0,0,2),
    @('Guam',32,0,0,31,0,0,1),
    @('Curazao',32,0,30,1,0,0,1),
    @('Fiyi',27,0,18,8,0,0,1),
    @('Timor Oriental',25,0,24,1,0,0,0),
    @('Santa Lucia',25,0,24,1,0,0,0),
    @('Granada',24,0,23,1,0,0,0),
    @('Nueva Caledonia',23,0,22,1,0,0,0),
    @('Laos',20,0,19,1,0,0,0),
    @('Dominica',18,0,18,0,0,0,0),
    @('Islas Virgenes de los Estados Unidos',17,0,0,17,0,0,0),
    @('San Cristobal y Nieves',17,0,17,0,0,0,0),
    @('Groenlandia',14,0,14,0,0,0,0),
    @('San Bartolome',13,0,6,7,0,0,0),
    @('Bonaire, San Eustaquio y Saba',13,0,7,6,0,0,0),
    @('Montserrat',13,0,12,0,0,0,1),
    @('Islas Malvinas',13,0,13,0,0,0,0),
    @('Santa Sede',12,0,12,0,0,0,0),
    @('Sahara Occidental',10,0,8,1,0,0,1),
    @('Islas Virgenes Britanicas',9,0,7,1,0,0,1),
    @('San Pedro y Miquelon',4,0,1,3,0,0,0),
    @('Anguila',3,0,3,0,0,0,0)
)

$row = 4
foreach ($r in $data) {
    for ($c = 0; $c -lt 8; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $r[$c]
    }
    $row++
}

Write-Host "Updated $($data.Count) country rows (rows 4-$($row - 1))"
